$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Pets table (A1:F11) -- refresh to the latest seed data. Headers (row 1)
#    are unchanged; rows 2-11 get new PetID/name/poster/species/sex/age data
#    (rows 10-11 are brand new rows).
# ---------------------------------------------------------------------------
$pets = @(
  @(1,  "Fido",       1, 1, "M", 1),
  @(2,  "Meowster",   2, 2, "M", 1),
  @(3,  "Ruffy",      2, 1, "M", 2),
  @(4,  "Scruffy",    3, 1, "M", 2),
  @(5,  "Purrington", 1, 2, "F", 3),
  @(6,  "Charles",    3, 2, "M", 3),
  @(7,  "Goldy",      1, 1, "F", 10),
  @(8,  "Khali",      3, 2, "F", 10),
  @(9,  "Mishka",     2, 1, "F", 5),
  @(10, "Woofy",      1, 1, "M", 5)
)

$r = 2
foreach ($pet in $pets) {
  $ws.Cells.Item($r, 1).Value = $pet[0]
  $ws.Cells.Item($r, 2).Value = $pet[1]
  $ws.Cells.Item($r, 3).Value = $pet[2]
  $ws.Cells.Item($r, 4).Value = $pet[3]
  $ws.Cells.Item($r, 5).Value = $pet[4]
  $ws.Cells.Item($r, 6).Value = $pet[5]
  $r = $r + 1
}

# New rows 10 & 11 need the same "note" style the rest of column G already
# carries (G2:G9 each hold an empty, Hyperlink-styled cell).
$ws.Range("G10").Style = "Hyperlink"
$ws.Range("G11").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 2. The old placeholder formatted-only cells at E19:E21 go away -- E19 is
#    fully cleared, E20/E21 get reused below for the new Users table.
# ---------------------------------------------------------------------------
$ws.Range("E19").Clear()

# ---------------------------------------------------------------------------
# 3. Categories table (A14:B16).
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = "CategoryID"
$ws.Cells.Item(14, 2).Value = "name"

$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Dog"

$ws.Cells.Item(16, 1).Value = 2
$ws.Cells.Item(16, 2).Value = "Cat"

# ---------------------------------------------------------------------------
# 4. Users/shelter table (A20:E23), including mailto hyperlinks on E21:E23.
# ---------------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = "UserID"
$ws.Cells.Item(20, 2).Value = "username"
$ws.Cells.Item(20, 3).Value = "password"
$ws.Cells.Item(20, 4).Value = "shelter"
$ws.Cells.Item(20, 5).Value = "email"
$ws.Range("E20").Style = "Normal"

$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "test"
$ws.Cells.Item(21, 3).Value = "abcd1234"
$ws.Cells.Item(21, 4).Value = "Safe Haven"
$ws.Cells.Item(21, 5).Value = "test@test.com"

$ws.Cells.Item(22, 1).Value = 2
$ws.Cells.Item(22, 2).Value = "jvasallo"
$ws.Cells.Item(22, 3).Value = "abcd1234"
$ws.Cells.Item(22, 4).Value = "First Woof"
$ws.Cells.Item(22, 5).Value = "jvasallo@test.com"

$ws.Cells.Item(23, 1).Value = 3
$ws.Cells.Item(23, 2).Value = "cyoung"
$ws.Cells.Item(23, 3).Value = "abcd1234"
$ws.Cells.Item(23, 4).Value = "North Pets"
$ws.Cells.Item(23, 5).Value = "cyoung@test.com"

# Hyperlinks are added bottom-up so the generated relationship ids come out
# rId1 -> E23, rId2 -> E22, rId3 -> E21 (matches the authored workbook).
$ws.Hyperlinks.Add($ws.Range("E23"), "mailto:cyoung@test.com")
$ws.Hyperlinks.Add($ws.Range("E22"), "mailto:jvasallo@test.com")
$ws.Hyperlinks.Add($ws.Range("E21"), "mailto:test@test.com")

$ws.Range("E21").Style = "Hyperlink"
$ws.Range("E22").Style = "Hyperlink"
$ws.Range("E23").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 5. Selection moves to A2:F3 with A2 active.
# ---------------------------------------------------------------------------
$ws.Range("A2:F3").Select()
